$wb = $excel.ActiveWorkbook

# --- Update the "Status" text used across the Overview, zh-cn and de-de sheets ---
# All four cells currently share the text "Ready for handoff"; they must all move to
# the new text "Handback transform failed" so the shared string is updated in place
# (no new shared-string entry should be created for this change).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handback transform failed"
$wsOverview.Range("F3").Value = "Handback transform failed"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Handback transform failed"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Handback transform failed"

# --- Populate the "Error Detail" column (P) for row 3 (the f64932bf file) with the
#     handback-transform error message, once per locale sheet ---
$wsZhCn.Range("P3").Value = "Handback file name: gscgiesn.4xi is different with handoff file name: f64932bf-35a8-433e-8407-e0da7a860823.4eb6bef3f85ace5215412e93d517234870bef313.zh-cn."
$wsDeDe.Range("P3").Value = "Handback file name: gscgiesn.4xi is different with handoff file name: f64932bf-35a8-433e-8407-e0da7a860823.4eb6bef3f85ace5215412e93d517234870bef313.de-de."

# --- Widen the "Error Detail" column (P, the 16th column) on both locale sheets now
#     that it holds long error messages ---
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
